$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the new record
$ws.Range("A2").Value = "ACC-1756905141601"
$ws.Range("B2").Value = "ali.arslan@agiemtech.ae"
$ws.Range("C2").Value = 500
$ws.Range("D2").Value = 400
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = "Active"

# Remove rows 3 and 4 entirely (data now ends at row 2)
$ws.Range("A3:F4").EntireRow.Delete()
